$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.766.75'
$ws.Range("E2").Value = '  -2.62%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.743.91'
$ws.Range("E3").Value = '  -5.06%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.77'
$ws.Range("E5").Value = '  -9.14%  '

# Row 6
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5033'
$ws.Range("E7").Value = '  -6.62%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.93'
$ws.Range("E8").Value = '  -6.48%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2679'
$ws.Range("E9").Value = '  -11.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06142'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.743.61'
$ws.Range("E11").Value = '  -5.17%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06928'
$ws.Range("E12").Value = '  -3.12%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.42'
$ws.Range("E13").Value = '  -12.83%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.520'
$ws.Range("E14").Value = '  -9.52%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5994'
$ws.Range("E15").Value = '  -18.68%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.84'
$ws.Range("E16").Value = '  -13.79%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.01%  '

# Row 18
$ws.Range("E18").Value = '  +0.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.768.99'
$ws.Range("E19").Value = '  -2.71%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006832'
$ws.Range("E20").Value = '  -13.44%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.56'
$ws.Range("E21").Value = '  -16.60%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.966.59'
$ws.Range("E22").Value = '  -5.27%  '

# Row 23
$ws.Range("E23").Value = '  -11.90%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.234'
$ws.Range("E24").Value = '  -12.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.134'
$ws.Range("E25").Value = '  -11.65%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.50'
$ws.Range("E26").Value = '  -3.83%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.515'
$ws.Range("E27").Value = '  -10.80%  '

# Row 28
$ws.Range("E28").Value = '  -11.84%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.808'
$ws.Range("E29").Value = '  -17.56%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '104.10'
$ws.Range("E30").Value = '  -6.11%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08111'
$ws.Range("E31").Value = '  -8.13%  '

# Row 32
$ws.Range("E32").Value = '  -11.69%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.464'
$ws.Range("E33").Value = '  -14.20%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04558'
$ws.Range("E34").Value = '  -5.47%  '

# Row 35
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.624'
$ws.Range("E36").Value = '  -10.11%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9800'
$ws.Range("E37").Value = '  -13.34%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6092'
$ws.Range("E38").Value = '  -16.45%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.662'
$ws.Range("E39").Value = '  -13.86%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01550'
$ws.Range("E40").Value = '  -9.66%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.914'
$ws.Range("E41").Value = '  -15.10%  '

# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.09%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.76'
$ws.Range("E43").Value = '  -5.73%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3811'
$ws.Range("E44").Value = '  -19.16%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.064'
$ws.Range("E45").Value = '  -14.07%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7323'
$ws.Range("E46").Value = '  -19.00%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05369'
$ws.Range("E47").Value = '  -6.94%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1110'
$ws.Range("E48").Value = '  -10.98%  '

# Row 49
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.931'
$ws.Range("E49").Value = '  -19.72%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.14'
$ws.Range("E50").Value = '  -13.44%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.44'
$ws.Range("E51").Value = '  -12.69%  '
